$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.108.80'
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("D3").Value = '1.830.66'
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4653'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2672'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06259'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.41%  '
$ws.Range("D10").Value = '1.829.70'
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07386'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.885'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6155'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.23%  '
$ws.Range("D16").Value = '30.069.02'
$ws.Range("E16").Value = '  -1.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '225.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007245'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.75%  '
$ws.Range("E20").Value = '  -6.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '2.069.33'
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("E23").Value = '  -8.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.826'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.149'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.857'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1008'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E30").Value = '  -1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.049'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.769'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04761'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.125'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7068'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.684'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01808'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.610'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8914'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.923'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.461'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3985'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.940'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.44%  '
$ws.Range("E46").Value = '  -6.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.369'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05513'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("E51").Value = '  -7.48%  '
